# Add "Aligner Options" (r2c_aligner_options) and "Extract Unmapped/Unassembled
# Reads" (extractUnmapped) rows to the assembly UI spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column E was widened (no longer auto "best fit") to better fit the new,
# longer label text.
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 32.16666666666666

# ---------------------------------------------------------------------------
# Column A ("input") for the two new rows.
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = "String? r2c_aligner_options"
$ws.Range("A27").Value = "Boolean extractUnmapped"

# ---------------------------------------------------------------------------
# Column D ("UI") marks both rows as visible ("yes").
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = "yes"
$ws.Range("D27").Value = "yes"

# ---------------------------------------------------------------------------
# Column E ("text") labels.
# ---------------------------------------------------------------------------
$ws.Range("E26").Value = "Aligner Options"
$ws.Range("E27").Value = "Extract Unmapped/Unassembled Reads"

# ---------------------------------------------------------------------------
# Column G ("tooltip") rich-text description of the aligner options, row 26.
# ---------------------------------------------------------------------------
$tooltip26 = $ws.Range("G26")
$tooltip26.Value = "Click Bowtie2 (https://bowtie-bio.sourceforge.net/bowtie2/manual.shtml#usage) | BWA mem(https://bio-bwa.sourceforge.net/bwa.shtml#3) | Minimap2(https://lh3.github.io/minimap2/minimap2.html) for detail. "

$run = $tooltip26.Characters(7, 71)
$run.Font.Name = "Arial"
$run.Font.Size = 14
$run.Font.Color = 0

$run = $tooltip26.Characters(78, 3)
$run.Font.Name = "Arial"
$run.Font.Size = 14
$run.Font.Color = 6710886

$run = $tooltip26.Characters(81, 52)
$run.Font.Name = "Arial"
$run.Font.Size = 14
$run.Font.Color = 0

$run = $tooltip26.Characters(133, 3)
$run.Font.Name = "Arial"
$run.Font.Size = 14
$run.Font.Color = 6710886

$run = $tooltip26.Characters(136, 54)
$run.Font.Name = "Arial"
$run.Font.Size = 14
$run.Font.Color = 0

$run = $tooltip26.Characters(190, 13)
$run.Font.Name = "Arial"
$run.Font.Size = 14
$run.Font.Color = 6710886

# Row 26 grows slightly taller to accommodate the 14pt tooltip text.
$ws.Rows(26).RowHeight = 18

# ---------------------------------------------------------------------------
# Column B ("default") - extractUnmapped defaults to False.
# ---------------------------------------------------------------------------
$ws.Range("B27").Value = $false

# ---------------------------------------------------------------------------
# Cosmetic: the previously styled (but otherwise blank) format on G8 is
# cleared back to the default style, and the active selection moves to F18.
# ---------------------------------------------------------------------------
$ws.Range("G8").Style = "Normal"

$ws.Range("F18").Select()
